$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with P1=14, Q1=15, copying the existing header style (O1) so that
# the cell format index matches exactly (bold font, border, centered/top aligned).
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update rows 2-25: swap I/K and M/O values, and add new P/Q columns with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I  (was 1)
    $ws.Cells.Item($r, 11).Value = 1  # K  (was 2)
    $ws.Cells.Item($r, 13).Value = 2  # M  (was 1)
    $ws.Cells.Item($r, 15).Value = 1  # O  (was 2)
    $ws.Cells.Item($r, 16).Value = 2  # P  (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q  (new)
}
